$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing C/D/E values for rows 8-15 (extr1..extr8)
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("C9").Value = 16

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $false

$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# Add new rows 16 and 17 (line7, line8)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

# Copy style from A15 to A16/A17 (bold/border/center style used for column A index)
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
